$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Users" to "1"
$ws.Name = "1"

# Change the existing row 2 values "a" -> "z" (reuses the shared string slot)
$ws.Range("A2").Value = "z"
$ws.Range("B2").Value = "z"

# Insert a new blank row above (old) row 2, pushing the "z" row down to row 3
$ws.Rows.Item(2).Insert()

# New row 2: "aa" / "aa"
$ws.Range("A2").Value = "aa"
$ws.Range("B2").Value = "aa"

# New row 4: "x" / "x"
$ws.Range("A4").Value = "x"
$ws.Range("B4").Value = "x"

# New row 5: "d" / "d"
$ws.Range("A5").Value = "d"
$ws.Range("B5").Value = "d"

# Update selection to B3
$ws.Range("B3").Select()
